{"js": "// Remove the \"Preparation\" homework prompts (the two numbered questions and\n// the trailing \"\\u2013>\" cue paragraph) that used to sit between the\n// \"Preparation\" heading and the \"Solutions\" heading. The paragraphs all used\n// to carry a dedicated list definition (numId 1001) that only existed for\n// them, so once the paragraphs are gone nothing in the document references\n// that list any more.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the \"Preparation\" paragraph (Author style) and the \"Solutions\"\n// paragraph (Heading 2 style) that bracket the block we need to remove.\n// Anchoring on those two stable paragraphs (rather than hard-coded indexes\n// or exact question wording) keeps this resilient to unrelated changes\n// elsewhere in the document.\nlet prepIdx = -1;\nlet solutionsIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (prepIdx === -1 && items[i].style === \"Author\" && text === \"Preparation\") {\n    prepIdx = i;\n    continue;\n  }\n  if (prepIdx !== -1 && items[i].style === \"Heading 2\" && text === \"Solutions\") {\n    solutionsIdx = i;\n    break;\n  }\n}\n\nif (prepIdx === -1 || solutionsIdx === -1) {\n  throw new Error(\n    \"Could not find the Preparation/Solutions boundary (prepIdx=\" + prepIdx + \", solutionsIdx=\" + solutionsIdx + \")\"\n  );\n}\n\n// Delete every paragraph strictly between \"Preparation\" and \"Solutions\",\n// walking backwards so earlier deletions don't shift later indexes.\nfor (let i = solutionsIdx - 1; i > prepIdx; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Preparation\" homework prompts (the two numbered questions and\n# the trailing \"-->\" cue paragraph) that used to sit between the\n# \"Preparation\" heading and the \"Solutions\" heading. The paragraphs used to\n# carry a dedicated list definition (numId 1001) that only existed for them,\n# so once the paragraphs are gone nothing in the document references that\n# list any more.\n\n$d = $word.ActiveDocument\n\n# Find the \"Preparation\" paragraph (Author style) and the \"Solutions\"\n# paragraph (Heading 2 style) that bracket the block we need to remove.\n# Anchoring on those two stable paragraphs (rather than hard-coded indexes\n# or exact question wording) keeps this resilient to unrelated changes\n# elsewhere in the document.\n$prepIdx = -1\n$solutionsIdx = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Range.Style.NameLocal\n    $text = $p.Range.Text.Trim()\n\n    if ($prepIdx -eq -1 -and $styleName -eq \"Author\" -and $text -eq \"Preparation\") {\n        $prepIdx = $i\n        continue\n    }\n\n    if ($prepIdx -ne -1 -and $styleName -eq \"Heading 2\" -and $text -eq \"Solutions\") {\n        $solutionsIdx = $i\n        break\n    }\n}\n\nif ($prepIdx -eq -1 -or $solutionsIdx -eq -1) {\n    Write-Output (\"Could not find the Preparation/Solutions boundary (prepIdx=\" + $prepIdx + \", solutionsIdx=\" + $solutionsIdx + \")\")\n} else {\n    # Delete every paragraph strictly between \"Preparation\" and \"Solutions\",\n    # walking backwards so earlier deletions don't shift later indexes.\n    for ($i = $solutionsIdx - 1; $i -gt $prepIdx; $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n"}
